$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 9 (Ano 2025) with the latest faturamento data
$ws.Range("B9").Value = 3845699.89
$ws.Range("C9").Value = 606140.24
$ws.Range("D9").Value = 4451840.13
$ws.Range("E9").Value = 13.61549881172395
$ws.Range("F9").Value = 86.38450118827605
$ws.Range("G9").Value = -41.4195681585722
$ws.Range("H9").Value = -30.55197566533356
$ws.Range("I9").Value = 38930
$ws.Range("J9").Value = 1663
$ws.Range("K9").Value = 40593
$ws.Range("L9").Value = 28079
$ws.Range("M9").Value = 158.5469614302503
$ws.Range("N9").Value = 8.24323981004107
